$d = $word.ActiveDocument

# Locate the "Product Development and Platform Architecture" paragraph
# within the Siege Analytics / PARTNER section.
$targetIdx = 0
for ($k = 1; $k -le $d.Paragraphs.Count; $k++) {
    $p = $d.Paragraphs.Item($k)
    if ($p.Range.Text -like "Product Development and Platform Architecture*") {
        $targetIdx = $k
        break
    }
}

if ($targetIdx -gt 0) {
    $newBullets = @(
        "• Conceived and architected redistricting platform incorporating boundary estimation algorithm used by 2,500+ analysts",
        "• Built multi-tenant data warehouse tracking decades of demographic data, enabling discovery of 500,000+ mischaracterized voters",
        "• Platform democratized redistricting analysis, reducing costs by 75% and enabling 200+ smaller organizations to participate"
    )

    $insertPos = $targetIdx
    foreach ($bulletText in $newBullets) {
        $anchorPara = $d.Paragraphs.Item($insertPos)
        $anchorRange = $anchorPara.Range
        $anchorRange.Collapse(0)
        $anchorRange.InsertParagraphAfter()

        $insertPos = $insertPos + 1
        $newPara = $d.Paragraphs.Item($insertPos)
        $newPara.Range.Text = $bulletText
    }
}
